$d = $word.ActiveDocument

# The "Requisitos" list is the last paragraph in the document (ListBullet style),
# immediately following the "Requisitos" Heading2 paragraph.
$paras = $d.Paragraphs
$listPara = $paras.Item($paras.Count)
$pos = $listPara.Range.Start

$newOrder = @(
  "LOB1053 -  Física III  (Requisito)",
  "LOQ4095 -  Química Geral Experimental  (Requisito)",
  "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
  "LOB1036 -  Geometria Analítica  (Requisito)",
  "LOB1024 -  Mecânica  (Requisito)",
  "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)",
  "LOB1004 -  Cálculo II  (Requisito)",
  "LOB1011 -  Eletricidade Aplicada  (Requisito)",
  "LOB1052 -  Cálculo III  (Requisito)",
  "LOB1012 -  Estatística  (Requisito)",
  "LOB1006 -  Cálculo IV  (Requisito)",
  "LOB1039 -  Física Experimental III  (Requisito)",
  "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)",
  "LOB1019 -  Física II  (Requisito)",
  "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)",
  "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
  "LOB1003 -  Cálculo I  (Requisito)",
  "LOB1038 -  Física Experimental I  (Requisito)",
  "LOB1018 -  Física I  (Requisito)",
  "LOB1037 -  Álgebra Linear  (Requisito)",
  "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)",
  "LOQ4103 -  Escrita Acadêmico Científica  (Requisito)",
  "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
  "LOQ4102 -  Nivelamento em Engenharia  (Requisito)"
)

$originalLengths = @(
  57,
  50,
  63,
  64,
  49,
  52,
  34,
  34,
  73,
  33,
  45,
  35,
  32,
  43,
  38,
  45,
  47,
  35,
  59,
  34,
  66,
  32,
  33,
  57
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $len = $originalLengths[$i]
    $rng = $d.Range($pos, $pos + $len)
    $rng.Text = $newOrder[$i]
    $pos = $rng.End + 1  # skip the line-break (vertical tab) character after this run
}

Write-Output "Reordered requisites list."
